$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 65, shifting existing rows 65-139 down to 66-140.
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with its data (same as old row 65's
# unchanged columns, plus the new Fecha / Precio minimo / Precio maximo /
# Precio promedio ponderado / Origen / Precio $/Kg values).
$ws.Range("A65").Value = 10
$ws.Range("B65").Value = "Vega Modelo de Temuco"
$ws.Range("C65").Value = "La Araucanía"
$ws.Range("D65").Value = 44994
$ws.Range("E65").Value = 9
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100101
$ws.Range("H65").Value = "Berries"
$ws.Range("I65").Value = 100101001
$ws.Range("J65").Value = "Arándano (blue)"
$ws.Range("K65").Value = "Sin especificar"
$ws.Range("L65").Value = "Primera"
$ws.Range("M65").Value = 200
$ws.Range("N65").Value = 1500
$ws.Range("O65").Value = 1500
$ws.Range("P65").Value = 1500
$ws.Range("Q65").Value = "$/kilo"
$ws.Range("R65").Value = "Región de La Araucanía"
$ws.Range("S65").Value = 1500
$ws.Range("T65").Value = 1

# Match the style used by the other date cells in column D (yyyy-mm-dd format).
$ws.Range("D65").NumberFormat = $ws.Range("D66").NumberFormat
